$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 16 with data for 12/10/2025
# Force the date-like text to be stored as text (not auto-converted to a date serial)
$ws.Cells.Item(16, 1).NumberFormat = "@"
$ws.Cells.Item(16, 1).Value = "12/10/2025"
$ws.Cells.Item(16, 1).Style = "Normal"

$ws.Cells.Item(16, 2).Value = 13219.05
$ws.Cells.Item(16, 3).Value = 0.1922042298597978
$ws.Cells.Item(16, 4).Value = 0.8077957701402022
$ws.Cells.Item(16, 5).Value = -101.77
$ws.Cells.Item(16, 6).Value = -23.4
$ws.Cells.Item(16, 7).Value = -19813.46
$ws.Cells.Item(16, 8).Value = -64.98
$ws.Cells.Item(16, 9).Value = -362.09
$ws.Cells.Item(16, 10).Value = -12.47
